# Update average price values on Sheet1 (recalculated "pages situation" fix)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("B2").Value = 1575.477723373003
$ws.Range("C2").Value = 1575.477723373003
$ws.Range("D2").Value = 1572.299293773217

# Row 3
$ws.Range("B3").Value = 1600.019014252174
$ws.Range("C3").Value = 1600.019014252174
$ws.Range("D3").Value = 1597.763145621847

# Row 4
$ws.Range("B4").Value = 1616.886424091637
$ws.Range("C4").Value = 1626.258828749242
$ws.Range("D4").Value = 1616.886424091637

# Row 5
$ws.Range("B5").Value = 1734.312784874813
$ws.Range("C5").Value = 1748.71612779832
$ws.Range("D5").Value = 1734.312784874813

# Row 6
$ws.Range("B6").Value = 1782.202981048386
$ws.Range("C6").Value = 1790.543296175052
$ws.Range("D6").Value = 1782.202981048386

# Row 7 (D7 stays an empty inline string cell, untouched)
$ws.Range("B7").Value = 1733.711422127711
$ws.Range("C7").Value = 1733.711422127711

# Row 8 (D8 stays an empty inline string cell, untouched)
$ws.Range("B8").Value = 1694.094743154283
$ws.Range("C8").Value = 1694.094743154283

# Row 9
$ws.Range("B9").Value = 1730.260844539496
$ws.Range("C9").Value = 1747.653453086408
$ws.Range("D9").Value = 1730.260844539496

$wb.Save()
